# Mise à jour de l'application : ajout de la colonne de présence du 08/01/2026 (CX)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: new date (serial 46030 = 2026-01-08), same style as CW1 ---
$ws.Cells.Item(1,102).Value = 46030
$ws.Cells.Item(1,101).Copy()
$ws.Cells.Item(1,102).PasteSpecial(-4122)

# --- Attendance codes for the new date, column CX, rows 2-30 ---
# (rows 12, 21 and 23 already stop short of column CW - no entry that day)
$ws.Cells.Item(2,102).Value = "P"
$ws.Cells.Item(2,101).Copy()
$ws.Cells.Item(2,102).PasteSpecial(-4122)
$ws.Cells.Item(3,102).Value = "P"
$ws.Cells.Item(3,101).Copy()
$ws.Cells.Item(3,102).PasteSpecial(-4122)
$ws.Cells.Item(4,102).Value = "P"
$ws.Cells.Item(4,101).Copy()
$ws.Cells.Item(4,102).PasteSpecial(-4122)
$ws.Cells.Item(5,102).Value = "P"
$ws.Cells.Item(5,101).Copy()
$ws.Cells.Item(5,102).PasteSpecial(-4122)
$ws.Cells.Item(6,102).Value = "R"
$ws.Cells.Item(6,101).Copy()
$ws.Cells.Item(6,102).PasteSpecial(-4122)
$ws.Cells.Item(7,102).Value = "P"
$ws.Cells.Item(7,101).Copy()
$ws.Cells.Item(7,102).PasteSpecial(-4122)
$ws.Cells.Item(8,102).Value = "R"
$ws.Cells.Item(8,101).Copy()
$ws.Cells.Item(8,102).PasteSpecial(-4122)
$ws.Cells.Item(9,102).Value = "RH"
$ws.Cells.Item(9,101).Copy()
$ws.Cells.Item(9,102).PasteSpecial(-4122)
$ws.Cells.Item(10,102).Value = "P"
$ws.Cells.Item(10,101).Copy()
$ws.Cells.Item(10,102).PasteSpecial(-4122)
$ws.Cells.Item(11,102).Value = "P"
$ws.Cells.Item(11,101).Copy()
$ws.Cells.Item(11,102).PasteSpecial(-4122)
$ws.Cells.Item(13,102).Value = "B"
$ws.Cells.Item(13,101).Copy()
$ws.Cells.Item(13,102).PasteSpecial(-4122)
$ws.Cells.Item(14,102).Value = "P"
$ws.Cells.Item(14,101).Copy()
$ws.Cells.Item(14,102).PasteSpecial(-4122)
$ws.Cells.Item(15,102).Value = "P"
$ws.Cells.Item(15,101).Copy()
$ws.Cells.Item(15,102).PasteSpecial(-4122)
$ws.Cells.Item(16,102).Value = "RH"
$ws.Cells.Item(16,101).Copy()
$ws.Cells.Item(16,102).PasteSpecial(-4122)
$ws.Cells.Item(17,102).Value = "P"
$ws.Cells.Item(17,101).Copy()
$ws.Cells.Item(17,102).PasteSpecial(-4122)
$ws.Cells.Item(18,102).Value = "P"
$ws.Cells.Item(18,101).Copy()
$ws.Cells.Item(18,102).PasteSpecial(-4122)
$ws.Cells.Item(19,102).Value = "P"
$ws.Cells.Item(19,101).Copy()
$ws.Cells.Item(19,102).PasteSpecial(-4122)
$ws.Cells.Item(20,102).Value = "P"
$ws.Cells.Item(20,101).Copy()
$ws.Cells.Item(20,102).PasteSpecial(-4122)
$ws.Cells.Item(22,102).Value = "P"
$ws.Cells.Item(22,101).Copy()
$ws.Cells.Item(22,102).PasteSpecial(-4122)
$ws.Cells.Item(24,102).Value = "P"
$ws.Cells.Item(24,101).Copy()
$ws.Cells.Item(24,102).PasteSpecial(-4122)
$ws.Cells.Item(25,102).Value = "R"
$ws.Cells.Item(25,101).Copy()
$ws.Cells.Item(25,102).PasteSpecial(-4122)
$ws.Cells.Item(26,102).Value = "P"
$ws.Cells.Item(26,101).Copy()
$ws.Cells.Item(26,102).PasteSpecial(-4122)
$ws.Cells.Item(27,102).Value = "P"
$ws.Cells.Item(27,101).Copy()
$ws.Cells.Item(27,102).PasteSpecial(-4122)
$ws.Cells.Item(28,102).Value = "P"
$ws.Cells.Item(28,101).Copy()
$ws.Cells.Item(28,102).PasteSpecial(-4122)
$ws.Cells.Item(29,102).Value = "P"
$ws.Cells.Item(29,101).Copy()
$ws.Cells.Item(29,102).PasteSpecial(-4122)
$ws.Cells.Item(30,102).Value = "P"
$ws.Cells.Item(30,101).Copy()
$ws.Cells.Item(30,102).PasteSpecial(-4122)

# --- Restore the active selection to the newly-entered cell ---
$ws.Range("CX28").Select() | Out-Null

$excel.CutCopyMode = 0
